$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# Per-language config: sheet name, handback datetime text, and the xlf
# hyperlink-target URL fragments (these differ only by the locale folder).
$langs = @(
    @{
        Sheet = "zh-cn"
        HandbackTime = "2016-03-11 14:10:35"
        XlfA = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e582c680502f3f42f19a439511de1171aede20ae/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/0967dc1b-4f2e-4de5-b6c9-5183586a9d0f.81ad96d38500c88f7c56bb4f43377ff8cc2313c4.zh-cn.xlf"
        XlfB = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e582c680502f3f42f19a439511de1171aede20ae/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/99dbd14f-f5ae-4da4-b25d-933915f5f79f.b6358d60477aa48c997366cf2d41cfdea25d4563.zh-cn.xlf"
        XlfADisplay = "0967dc1b-4f2e-4de5-b6c9-5183586a9d0f.81ad96d38500c88f7c56bb4f43377ff8cc2313c4.zh-cn.xlf"
        XlfBDisplay = "99dbd14f-f5ae-4da4-b25d-933915f5f79f.b6358d60477aa48c997366cf2d41cfdea25d4563.zh-cn.xlf"
    },
    @{
        Sheet = "de-de"
        HandbackTime = "2016-03-11 14:10:43"
        XlfA = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6b0ce07da92f91fabf3606670c6e52c0527ec333/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/0967dc1b-4f2e-4de5-b6c9-5183586a9d0f.81ad96d38500c88f7c56bb4f43377ff8cc2313c4.de-de.xlf"
        XlfB = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6b0ce07da92f91fabf3606670c6e52c0527ec333/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/99dbd14f-f5ae-4da4-b25d-933915f5f79f.b6358d60477aa48c997366cf2d41cfdea25d4563.de-de.xlf"
        XlfADisplay = "0967dc1b-4f2e-4de5-b6c9-5183586a9d0f.81ad96d38500c88f7c56bb4f43377ff8cc2313c4.de-de.xlf"
        XlfBDisplay = "99dbd14f-f5ae-4da4-b25d-933915f5f79f.b6358d60477aa48c997366cf2d41cfdea25d4563.de-de.xlf"
    }
)

# The two source-file rows are identical between sheets.
$mdUrlA = "https://github.com/OpenLocalizationTest/oltest/blob/7ae628596a7f97083478e4a994af978c15c090b0/e2e/0967dc1b-4f2e-4de5-b6c9-5183586a9d0f.md"
$mdUrlB = "https://github.com/OpenLocalizationTest/oltest/blob/7ae628596a7f97083478e4a994af978c15c090b0/e2e/99dbd14f-f5ae-4da4-b25d-933915f5f79f.md"
$mdDisplayA = "0967dc1b-4f2e-4de5-b6c9-5183586a9d0f.md"
$mdDisplayB = "99dbd14f-f5ae-4da4-b25d-933915f5f79f.md"

foreach ($lang in $langs) {
    $ws = $wb.Worksheets.Item($lang.Sheet)

    # Status column now reflects the completed handback.
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus

    # New "Latest Target File" (F) and "Latest Handback File" (G) columns,
    # populated with hyperlinks that mirror the existing Source File Name /
    # Latest Handoff File links for each row.
    $ws.Hyperlinks.Add($ws.Range("F2"), $mdUrlA, "", "", $mdDisplayA)
    $ws.Hyperlinks.Add($ws.Range("G2"), $lang.XlfA, "", "", $lang.XlfADisplay)

    $ws.Hyperlinks.Add($ws.Range("F3"), $mdUrlB, "", "", $mdDisplayB)
    $ws.Hyperlinks.Add($ws.Range("G3"), $lang.XlfB, "", "", $lang.XlfBDisplay)

    # Latest Handback DateTime (H) now has a real timestamp instead of the
    # zero-date placeholder.
    $ws.Range("H2").Value = $lang.HandbackTime
    $ws.Range("H3").Value = $lang.HandbackTime
}
